# Auto commit at 2025-09-29 10:08:52.78
# Updates the monthly metric actuals on the "Metrics" sheet; the "today"
# sheet pulls these via formulas (Metrics!B2:B13) and recalculates on its
# own. Also clears a stray manual value in today!I13 and restores the
# selections that were active on each sheet when the workbook was saved.

$wb = $excel.ActiveWorkbook

# --- Metrics sheet: refresh the metric values -----------------------------
$metrics = $wb.Worksheets.Item("Metrics")

$metrics.Range("B2").Value  = 417195.11
$metrics.Range("B3").Value  = 337231.43000000005
$metrics.Range("B4").Value  = 131626.99
$metrics.Range("B5").Value  = 16583
$metrics.Range("B6").Value  = 4336445.9899999993
$metrics.Range("B7").Value  = 3664758.9099999997
$metrics.Range("B8").Value  = 1260992.6700000002
$metrics.Range("B9").Value  = 167743
$metrics.Range("B10").Value = 32801769.79099983
$metrics.Range("B11").Value = 19694628.980000004
$metrics.Range("B12").Value = 11542701.559999999
$metrics.Range("B13").Value = 1265370

$metrics.Activate() | Out-Null
$metrics.Range("D8").Select() | Out-Null

# --- today sheet: clear stray manual entry and restore selection ----------
$today = $wb.Worksheets.Item("today")

$today.Range("I13").ClearContents() | Out-Null

$today.Activate() | Out-Null
$today.Range("F5").Select() | Out-Null
